$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 271-306 (columns H/I mostly) per new daily stats
$ws.Cells.Item(271, 8).Value = 42580
$ws.Cells.Item(271, 9).Value = 1620
$ws.Cells.Item(272, 8).Value = 30999
$ws.Cells.Item(272, 9).Value = 1666
$ws.Cells.Item(273, 8).Value = 27150
$ws.Cells.Item(273, 9).Value = 1367
$ws.Cells.Item(274, 8).Value = 28390
$ws.Cells.Item(274, 9).Value = 1348
$ws.Cells.Item(275, 8).Value = 28740
$ws.Cells.Item(275, 9).Value = 1239
$ws.Cells.Item(278, 8).Value = 29996
$ws.Cells.Item(278, 9).Value = 2101
$ws.Cells.Item(279, 8).Value = 43398
$ws.Cells.Item(279, 9).Value = 3143
$ws.Cells.Item(280, 8).Value = 35582
$ws.Cells.Item(280, 9).Value = 2405
$ws.Cells.Item(281, 8).Value = 45218
$ws.Cells.Item(281, 9).Value = 3255
$ws.Cells.Item(282, 8).Value = 46798
$ws.Cells.Item(282, 9).Value = 2857
$ws.Cells.Item(283, 8).Value = 17556
$ws.Cells.Item(285, 8).Value = 40904
$ws.Cells.Item(285, 9).Value = 3432
$ws.Cells.Item(286, 8).Value = 54313
$ws.Cells.Item(286, 9).Value = 4255
$ws.Cells.Item(287, 8).Value = 57739
$ws.Cells.Item(287, 9).Value = 3927
$ws.Cells.Item(288, 8).Value = 56403
$ws.Cells.Item(288, 9).Value = 3978
$ws.Cells.Item(289, 8).Value = 64626
$ws.Cells.Item(289, 9).Value = 3700
$ws.Cells.Item(290, 8).Value = 17876
$ws.Cells.Item(290, 9).Value = 1489
$ws.Cells.Item(292, 8).Value = 81520
$ws.Cells.Item(292, 9).Value = 7220
$ws.Cells.Item(293, 8).Value = 82467
$ws.Cells.Item(293, 9).Value = 5846
$ws.Cells.Item(294, 8).Value = 91480
$ws.Cells.Item(294, 9).Value = 5081
$ws.Cells.Item(295, 8).Value = 19352
$ws.Cells.Item(297, 8).Value = 2278
$ws.Cells.Item(297, 9).Value = 206
$ws.Cells.Item(298, 8).Value = 3047
$ws.Cells.Item(298, 9).Value = 280
$ws.Cells.Item(299, 8).Value = 63663
$ws.Cells.Item(299, 9).Value = 6686
$ws.Cells.Item(300, 8).Value = 70054
$ws.Cells.Item(300, 9).Value = 6844
$ws.Cells.Item(301, 8).Value = 68737
$ws.Cells.Item(301, 9).Value = 5426
$ws.Cells.Item(302, 8).Value = 71032
$ws.Cells.Item(302, 9).Value = 5158
$ws.Cells.Item(303, 8).Value = 11468
$ws.Cells.Item(303, 9).Value = 753
$ws.Cells.Item(304, 8).Value = 6506
$ws.Cells.Item(304, 9).Value = 476
$ws.Cells.Item(305, 8).Value = 3077
$ws.Cells.Item(305, 9).Value = 296
$ws.Cells.Item(306, 8).Value = 68407
$ws.Cells.Item(306, 9).Value = 6941
# Append new rows 307 and 308 for 2021-01-05 and 2021-01-06
$ws.Cells.Item(307, 1).Value = 44201
$ws.Cells.Item(307, 2).Value = 196047
$ws.Cells.Item(307, 3).Value = 135760
$ws.Cells.Item(307, 4).Value = 57630
$ws.Cells.Item(307, 5).Value = 17822
$ws.Cells.Item(307, 6).Value = 4959
$ws.Cells.Item(307, 7).Value = 2657
$ws.Cells.Item(307, 8).Value = 69321
$ws.Cells.Item(307, 9).Value = 5975
$ws.Cells.Item(308, 1).Value = 44202
$ws.Cells.Item(308, 2).Value = 198184
$ws.Cells.Item(308, 3).Value = 138383
$ws.Cells.Item(308, 4).Value = 57084
$ws.Cells.Item(308, 5).Value = 8257
$ws.Cells.Item(308, 6).Value = 2137
$ws.Cells.Item(308, 7).Value = 2717
$ws.Cells.Item(308, 8).Value = 17028
$ws.Cells.Item(308, 9).Value = 1424